{"js": "const body = context.document.body;\n\n// --- Change 1 -----------------------------------------------------------\n// The single run \"uerimientos Funcio\" is split in two: \"uerimientos \"\n// stays as plain text (now needs xml:space=\"preserve\" since it ends\n// with a space) and \"Funcio\" gets wrapped in a spell-check proofing\n// error (w:proofErr spellStart/spellEnd), as Word does after a\n// proofing pass flags \"Funcio\" as a possible misspelling.\nconst target = body.search(\"uerimientos Funcio\", { matchCase: true });\ntarget.load(\"text\");\nawait context.sync();\n\nif (target.items.length > 0) {\n  const flatOpcXml =\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body><w:p>' +\n    '<w:r><w:t xml:space=\"preserve\">uerimientos </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>Funcio</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '</w:p></w:body></w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>';\n  target.items[0].insertOoxml(flatOpcXml, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- Change 2 -----------------------------------------------------------\n// Remove the run containing the stray word \"hola\" from the bulleted\n// paragraph; the (now empty) paragraph itself, with its list\n// formatting, is left in place.\nconst holaRanges = body.search(\"hola\", { matchCase: true });\nholaRanges.load(\"text\");\nawait context.sync();\n\nif (holaRanges.items.length > 0) {\n  holaRanges.items[0].insertText(\"\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# --- Change 1 -----------------------------------------------------------\n# The single run \"uerimientos Funcio\" is split in two: \"uerimientos \"\n# stays as-is (now with xml:space=\"preserve\" because it ends with a\n# space) and \"Funcio\" is wrapped in a spell-check proofing error\n# (w:proofErr spellStart/spellEnd), as Word does after a proofing pass\n# flags \"Funcio\" as a possible misspelling.\n$find1 = $d.Content.Find\n$find1.Text = \"uerimientos Funcio\"\n$find1.Execute() | Out-Null\nif ($find1.Found) {\n    $rng1 = $find1.Parent\n    $xml1 = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:t xml:space=\"preserve\">uerimientos </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>Funcio</w:t></w:r><w:proofErr w:type=\"spellEnd\"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n    $rng1.InsertXML($xml1)\n}\n\n# --- Change 2 -----------------------------------------------------------\n# Remove the run containing the stray word \"hola\" from the bulleted\n# paragraph; the (now empty) paragraph itself, with its list formatting,\n# is left in place.\n$find2 = $d.Content.Find\n$find2.Text = \"hola\"\n$find2.Replacement.Text = \"\"\n$find2.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, \"\", 2) | Out-Null\n\nWrite-Output \"done\"\n"}
